$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column R so that the existing R..AE block
# shifts right to U..AH, making room for the new subject columns.
$ws.Range("R1:T1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data values for row 2 in the inserted columns.
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Update existing text values in row 2.
$ws.Range("D2").Value = "unknown"
$ws.Range("E2").Value = "unknown"
$ws.Range("F2").Value = "unknown"
$ws.Range("G2").Value = "considered"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "unknown"
$ws.Range("J2").Value = "unknown"

$wb.Save()
